# Fruta / hortaliza, semanal
# Insert 3 new weekly rows of "Espárragos" price data at the top of the
# data block (row 71), pushing the existing rows 71-88 down to 74-91.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 71:88 down to 74:91, inheriting formatting from the
# row above (same behaviour as Excel's native "Insert Rows").
$ws.Rows("71:73").Insert()

# New row 71 - Banquete
$ws.Cells.Item(71, 1).Value = 6
$ws.Cells.Item(71, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(71, 3).Value = "Metropolitana"
$ws.Cells.Item(71, 4).Value = 44505
$ws.Cells.Item(71, 5).Value = 13
$ws.Cells.Item(71, 6).Value = 300000000
$ws.Cells.Item(71, 7).Value = "Espárragos"
$ws.Cells.Item(71, 8).Value = "Sin especificar"
$ws.Cells.Item(71, 9).Value = "Banquete"
$ws.Cells.Item(71, 10).Value = 550
$ws.Cells.Item(71, 11).Value = 1300
$ws.Cells.Item(71, 12).Value = 1400
$ws.Cells.Item(71, 13).Value = 1358
$ws.Cells.Item(71, 14).Value = "`$/kilo"
$ws.Cells.Item(71, 15).Value = "Región Metropolitana"
$ws.Cells.Item(71, 16).Value = 1358
$ws.Cells.Item(71, 17).Value = 1
$ws.Cells.Item(71, 18).Value = "Hortaliza"

# New row 72 - Primera
$ws.Cells.Item(72, 1).Value = 6
$ws.Cells.Item(72, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(72, 3).Value = "Metropolitana"
$ws.Cells.Item(72, 4).Value = 44505
$ws.Cells.Item(72, 5).Value = 13
$ws.Cells.Item(72, 6).Value = 300000000
$ws.Cells.Item(72, 7).Value = "Espárragos"
$ws.Cells.Item(72, 8).Value = "Sin especificar"
$ws.Cells.Item(72, 9).Value = "Primera"
$ws.Cells.Item(72, 10).Value = 680
$ws.Cells.Item(72, 11).Value = 1100
$ws.Cells.Item(72, 12).Value = 1200
$ws.Cells.Item(72, 13).Value = 1163
$ws.Cells.Item(72, 14).Value = "`$/kilo"
$ws.Cells.Item(72, 15).Value = "Región Metropolitana"
$ws.Cells.Item(72, 16).Value = 1163
$ws.Cells.Item(72, 17).Value = 1
$ws.Cells.Item(72, 18).Value = "Hortaliza"

# New row 73 - Segunda
$ws.Cells.Item(73, 1).Value = 6
$ws.Cells.Item(73, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(73, 3).Value = "Metropolitana"
$ws.Cells.Item(73, 4).Value = 44505
$ws.Cells.Item(73, 5).Value = 13
$ws.Cells.Item(73, 6).Value = 300000000
$ws.Cells.Item(73, 7).Value = "Espárragos"
$ws.Cells.Item(73, 8).Value = "Sin especificar"
$ws.Cells.Item(73, 9).Value = "Segunda"
$ws.Cells.Item(73, 10).Value = 430
$ws.Cells.Item(73, 11).Value = 900
$ws.Cells.Item(73, 12).Value = 1000
$ws.Cells.Item(73, 13).Value = 965
$ws.Cells.Item(73, 14).Value = "`$/kilo"
$ws.Cells.Item(73, 15).Value = "Región Metropolitana"
$ws.Cells.Item(73, 16).Value = 965
$ws.Cells.Item(73, 17).Value = 1
$ws.Cells.Item(73, 18).Value = "Hortaliza"
